# Apply the commit's changes to the workbook.
$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: top-of-page properties table ---
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/capitation-arrangement"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet: detail table ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element row; its Constraint(s) column (AI) is cleared.
$elements.Range("AI2").Value = ""
